$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new profit row for 2025-10-16 (row 60), following the same
# shape as the existing rows: column A holds the date as text, column B
# holds the numeric profit value.

# Force column A's new cell to be treated as text so the date-like string
# "10/16/2025" is stored literally instead of being auto-converted into a
# serial date number (matches how the rest of the Date column is stored).
$ws.Range("A60").NumberFormat = "@"
$ws.Range("A60").Value = "10/16/2025"

$ws.Range("B60").Value = 11051.71
